$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("units")

# Set A2 = 0 (new cell, was previously empty)
$ws.Range("A2").Value = 0

# Update the selected/active cell to B13 (matches the diff's selection change)
$ws.Range("B13").Select()
